# Apply updated crypto price / volume figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference an untouched, default-styled cell so assignments below can be
# reset to the plain (no explicit style index) look the sheet already uses,
# instead of leaving behind a bespoke "quote prefix" style per edited cell.
$plainStyle = $ws.Range("B2").Style

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    # Leading apostrophe forces Excel to store the literal text instead of
    # auto-coercing number-looking strings (e.g. "134.26") into a float.
    $rng.Value = "'" + $text
    $rng.Style = $plainStyle
}

Set-TextValue "D2" "64.296.74"
Set-TextValue "E2" "  +0.22%  "
Set-TextValue "D3" "3.501.06"
Set-TextValue "E3" "  -0.59%  "
Set-TextValue "E4" "  +0.03%  "
Set-TextValue "D5" "590.07"
Set-TextValue "E5" "  +0.79%  "
Set-TextValue "D6" "134.26"
Set-TextValue "E6" "  -0.23%  "
Set-TextValue "E7" "  -0.01%  "
Set-TextValue "E8" "  -0.61%  "
Set-TextValue "E9" "  +0.06%  "
Set-TextValue "E10" "  +2.63%  "
Set-TextValue "E11" "  +2.41%  "
Set-TextValue "D12" "4.100.57"
Set-TextValue "E12" "  -0.54%  "
Set-TextValue "E13" "  +1.12%  "
Set-TextValue "E14" "  +1.02%  "
Set-TextValue "D15" "3.503.30"
Set-TextValue "E15" "  -0.51%  "
Set-TextValue "D16" "64.309.80"
Set-TextValue "E16" "  +0.18%  "
Set-TextValue "D17" "25.64"
Set-TextValue "E17" "  -6.72%  "
Set-TextValue "D18" "9.86"
Set-TextValue "E18" "  +0.81%  "
Set-TextValue "E19" "  +2.44%  "
Set-TextValue "D20" "13.51"
Set-TextValue "E20" "  -2.79%  "
Set-TextValue "D21" "392.73"
Set-TextValue "E21" "  +2.58%  "
Set-TextValue "E22" "  +0.20%  "
Set-TextValue "D23" "3.641.98"
Set-TextValue "E23" "  -0.64%  "
Set-TextValue "D24" "74.56"
Set-TextValue "E24" "  +0.70%  "
Set-TextValue "D25" "0.997"
Set-TextValue "E25" "  -0.29%  "
Set-TextValue "E26" "  -0.04%  "
Set-TextValue "E27" "  -0.03%  "
Set-TextValue "D28" "7.36"
Set-TextValue "E28" "  -1.38%  "
Set-TextValue "E29" "  +1.28%  "
Set-TextValue "D30" "8.22"
Set-TextValue "E30" "  -2.64%  "
Set-TextValue "D31" "1.47"
Set-TextValue "E31" "  -7.23%  "
Set-TextValue "D32" "3.523.82"
Set-TextValue "E32" "  -0.34%  "
Set-TextValue "E33" "  +5.46%  "
Set-TextValue "D35" "23.47"
Set-TextValue "D36" "5.13"
Set-TextValue "E36" "  -5.43%  "
Set-TextValue "D37" "6.88"
Set-TextValue "E37" "  -1.03%  "
Set-TextValue "E38" "  -0.71%  "
Set-TextValue "D39" "167.50"
Set-TextValue "D40" "0.0779"
Set-TextValue "E40" "  -0.86%  "
Set-TextValue "D41" "0.810"
Set-TextValue "E41" "  -0.24%  "
Set-TextValue "E42" "  +0.04%  "
Set-TextValue "D43" "25.01"
Set-TextValue "E43" "  -6.44%  "
Set-TextValue "D44" "4.39"
Set-TextValue "E44" "  -0.26%  "
Set-TextValue "E45" "  +3.15%  "
Set-TextValue "E46" "  -3.76%  "
Set-TextValue "D47" "6.75"
Set-TextValue "D48" "2.363.58"
Set-TextValue "E48" "  -4.87%  "
Set-TextValue "D49" "0.892"
Set-TextValue "E49" "  -1.81%  "
Set-TextValue "E50" "  -1.46%  "
Set-TextValue "D51" "21.11"
Set-TextValue "E51" "  -1.46%  "
